$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a duplicate "Variable" column (F) for the first block of questions (rows 3-6) ---
$ws.Range("F3").Value = "WASI_Complete"
$ws.Range("F4").Value = "WASI_incomplete_reason"
$ws.Range("F5").Value = "WASI_Valid"
$ws.Range("F6").Value = "WASI_Invalid_Reason"

# --- Re-purpose rows 20-23 to describe Block Design + Full Scale IQ, shifting the
#     previous "Percentile" helper rows out ---
$ws.Range("A20").Value = "Block Design Raw Score"
$ws.Range("B20").Value = "WASI_BD_Raw"

$ws.Range("A21").Value = "Block Design T Score"
$ws.Range("B21").Value = "WASI_BD_T"
$ws.Range("C21").Value = "decimal"

$ws.Range("A22").Value = "Full Scale IQ"
$ws.Range("B22").Value = "WASI_FSIQ"

$ws.Range("A23").Value = "Full Scale IQ Percentile"
$ws.Range("B23").Value = "WASI_FSIQ_P"
$ws.Range("C23").Value = "decimal"

# --- Rows 24-26 no longer hold data (old rows 24-27 content is gone); clear them ---
$ws.Range("A24:C26").ClearContents()

# --- Row 27 is removed entirely ---
$ws.Rows.Item(27).Delete()

# --- Update the selection to match the new authoring session ---
$ws.Range("A8:B23").Select()
